$d = $word.ActiveDocument

$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2B643301" w14:textId="337D1DED" w:rsidR="00174674" w:rsidRDefault="00174674" w:rsidP="00174674"><w:pPr><w:rPr><w:lang w:val="cs-CZ"/></w:rPr></w:pPr><w:r w:rsidRPr="00C9316E"><w:rPr><w:b/><w:bCs/><w:lang w:val="cs-CZ"/></w:rPr><w:t>Hlavní myšlenka:</w:t></w:r><w:r w:rsidR="00C9316E"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C9316E" w:rsidRPr="00C9316E"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>vyjádření vlastních pocitů (značné autobiografické prvky)</w:t></w:r><w:r w:rsidR="00183513"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>, obava ze smrti</w:t></w:r></w:p>')

$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3E056A6E" w14:textId="0C9882B6" w:rsidR="00174674" w:rsidRDefault="00174674" w:rsidP="0083780A"><w:pPr><w:rPr><w:lang w:val="cs-CZ"/></w:rPr></w:pPr><w:r w:rsidRPr="00490CF7"><w:rPr><w:b/><w:bCs/><w:lang w:val="cs-CZ"/></w:rPr><w:t>Jazykové prostředky:</w:t></w:r><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve">převažuje spisovný jazyk, ovšem s mnoha zastaralými </w:t></w:r><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>výrazy – archaismy</w:t></w:r><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>* (kuropění [svítání],</w:t></w:r><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>hoře [bolest srdce/zármutek]) nebo přechodníky (maje, rozjímaje, hledaje); místy i hovorové (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>tos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>řek) a knižní výrazy (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>dí</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> [říká], světice [žena prohlášená za svatou], nadvzdušná [jemnější než</w:t></w:r><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve">vzduch/nehmotná]) nebo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>cizojaz</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0083780A" w:rsidRPr="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>. prvky z hebrejštiny (mene tekel [naléhavé varování]);</w:t></w:r><w:r w:rsidR="00490CF7"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> využívá se </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00490CF7"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>onomatoipe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00490CF7"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="00490CF7" w:rsidRPr="00490CF7"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>je zvukomalebné slovo, tj. slovo foneticky napodobující přirozené zvuky</w:t></w:r><w:r w:rsidR="00490CF7"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>)</w:t></w:r></w:p>')

$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0CF0CB1A" w14:textId="544C4DAE" w:rsidR="00E448AD" w:rsidRDefault="00E448AD" w:rsidP="00174674"><w:pPr><w:rPr><w:lang w:val="cs-CZ"/></w:rPr></w:pPr><w:r w:rsidRPr="0083780A"><w:rPr><w:b/><w:bCs/><w:lang w:val="cs-CZ"/></w:rPr><w:t>Kompozice:</w:t></w:r><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> pouze 108 veršů (rozděleno do 18 slok po 6 verších), vypravěč je v </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>ich</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0083780A"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>-formě je nespolehlivý</w:t></w:r><w:r w:rsidR="00490CF7"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> (je mezi snem a bděním) přímí subjekt</w:t></w:r><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve">, osmistopá trochej, a b c b </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> první a třetí verš má vnitřní rým (v polovině verše je slovo které se rýmuje se slovem na konci verše), v čtvrté a pátém verši je vždy epifora (verš je zakončený stejně), šestý veš je vždy </w:t></w:r><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>poloviční,</w:t></w:r><w:r w:rsidR="000406AF"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> takže jen čtyřstopý obsahuje refrén nebo jeho variaci</w:t></w:r></w:p>')

$p18 = $d.Paragraphs.Item(18)
$p18.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="751A045D" w14:textId="662F6447" w:rsidR="006375DA" w:rsidRDefault="00183513" w:rsidP="006375DA"><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="cs-CZ"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>A</w:t></w:r><w:r w:rsidR="00A907EB"><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>m</w:t></w:r><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve">erický </w:t></w:r><w:bookmarkStart w:id="0" w:name="_Hlk135484785"/><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve">romantický básník, </w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>prozaik a esejista, pokládá se za zakladatele moderního hororového a detektivního žánru a za raného představitele sci-fi</w:t></w:r></w:p>')

$p19 = $d.Paragraphs.Item(19)
$p19.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="47EFCD33" w14:textId="53285235" w:rsidR="00183513" w:rsidRDefault="00183513" w:rsidP="006375DA"><w:pPr><w:pStyle w:val="Odstavecseseznamem"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="cs-CZ"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Narodil se do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>bostnu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>, rodiče mu zemřeli a dostal se do adoptivní rodiny, studoval, měl potíže s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t>drogramy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="cs-CZ"/></w:rPr><w:t xml:space="preserve"> a depresemi a trpěl alkoholismem během nichž literárně tvořil</w:t></w:r></w:p>')
